$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '51.057.48'
$ws.Range('E2').Value = '  -0.99%  '

# Row 3
$ws.Range('D3').Value = '2.936.26'
$ws.Range('E3').Value = '  -1.55%  '

# Row 4
$ws.Range('E4').Value = '  +0.04%  '

# Row 5
$ws.Range('D5').Value = '''374.04'
$ws.Range('E5').Value = '  -1.40%  '

# Row 6
$ws.Range('D6').Value = '''100.74'
$ws.Range('E6').Value = '  -3.93%  '

# Row 7
$ws.Range('D7').Value = '''0.534'
$ws.Range('E7').Value = '  -1.23%  '

# Row 9
$ws.Range('D9').Value = '''0.583'
$ws.Range('E9').Value = '  -2.02%  '

# Row 10
$ws.Range('D10').Value = '''36.13'
$ws.Range('E10').Value = '  -3.03%  '

# Row 11
$ws.Range('E11').Value = '  -0.71%  '

# Row 12
$ws.Range('E12').Value = '  -0.08%  '

# Row 13
$ws.Range('D13').Value = '3.404.74'
$ws.Range('E13').Value = '  -1.10%  '

# Row 14
$ws.Range('D14').Value = '''17.94'
$ws.Range('E14').Value = '  -2.66%  '

# Row 15
$ws.Range('D15').Value = '''7.49'
$ws.Range('E15').Value = '  -1.32%  '

# Row 16
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.943.05'
$ws.Range('E16').Value = '  -1.44%  '

# Row 17
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').Value = '''11.04'
$ws.Range('E17').Value = '  +48.43%  '

# Row 18
$ws.Range('D18').Value = '''0.971'
$ws.Range('E18').Value = '  +0.21%  '

# Row 19
$ws.Range('D19').Value = '50.987.84'
$ws.Range('E19').Value = '  -1.01%  '

# Row 20
$ws.Range('D20').Value = '''3.15'
$ws.Range('E20').Value = '  -5.39%  '

# Row 21
$ws.Range('D21').Value = '''12.37'
$ws.Range('E21').Value = '  -4.67%  '

# Row 22
$ws.Range('D22').Value = '0.0₃0955'
$ws.Range('E22').Value = '  -0.84%  '

# Row 23
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').Value = '''264.44'
$ws.Range('E23').Value = '  +0.80%  '

# Row 24
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '''68.58'
$ws.Range('E24').Value = '  -1.16%  '

# Row 25
$ws.Range('D25').Value = '''3.12'
$ws.Range('E25').Value = '  +9.98%  '

# Row 26
$ws.Range('E26').Value = '  -1.49%  '

# Row 27
$ws.Range('D27').Value = '''7.40'
$ws.Range('E27').Value = '  -4.00%  '

# Row 28
$ws.Range('E28').Value = '  -0.04%  '

# Row 29
$ws.Range('E29').Value = '  -4.11%  '

# Row 30
$ws.Range('D30').Value = '''25.60'
$ws.Range('E30').Value = '  -1.26%  '

# Row 31
$ws.Range('E31').Value = '  -4.29%  '

# Row 32
$ws.Range('D32').Value = '''9.92'
$ws.Range('E32').Value = '  +0.71%  '

# Row 33
$ws.Range('E33').Value = '  -0.89%  '

# Row 34
$ws.Range('E34').Value = '  -1.82%  '

# Row 35
$ws.Range('D35').Value = '''33.14'
$ws.Range('E35').Value = '  -5.56%  '

# Row 36
$ws.Range('D36').Value = '''0.0442'
$ws.Range('E36').Value = '  -0.79%  '

# Row 37
$ws.Range('E37').Value = '  -0.13%  '

# Row 38
$ws.Range('D38').Value = '''3.13'
$ws.Range('E38').Value = '  +3.22%  '

# Row 39
$ws.Range('E39').Value = '  -0.24%  '

# Row 40
$ws.Range('D40').Value = '''16.32'
$ws.Range('E40').Value = '  -5.17%  '

# Row 41
$ws.Range('E41').Value = '  -3.57%  '

# Row 42
$ws.Range('D42').Value = '''2.47'
$ws.Range('E42').Value = '  -4.52%  '

# Row 43
$ws.Range('D43').Value = '''119.64'
$ws.Range('E43').Value = '  -4.69%  '

# Row 44
$ws.Range('D44').Value = '''21.06'
$ws.Range('E44').Value = '  -3.11%  '

# Row 45
$ws.Range('D45').Value = '''0.273'
$ws.Range('E45').Value = '  -3.78%  '

# Row 46
$ws.Range('E46').Value = '  -1.01%  '

# Row 47
$ws.Range('E47').Value = '  +1.96%  '

# Row 48
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '1.986.97'
$ws.Range('E48').Value = '  -2.41%  '

# Row 49
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').Value = '''2.28'
$ws.Range('E49').Value = '  -4.47%  '

# Row 50
$ws.Range('D50').Value = '''0.0326'
$ws.Range('E50').Value = '  -2.95%  '

# Row 51
$ws.Range('E51').Value = '  +0.48%  '
